$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1702
$ws.Range("E2").Value = 46200502396
$ws.Range("X2").Value = "DN4127450128907"

$ws.Range("A3").Value = 1703
$ws.Range("E3").Value = 46200502397
$ws.Range("X3").Value = "DN4127450128908"
